$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Fix typos in the opening ("First Paragraph") sentence of the Intro.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("A pimary role", $true, $false, $false, $false, $false, $true, 1, $false, "A primary role", 2)
$d.Content.Find.Execute("adequatly manage", $true, $false, $false, $false, $false, $true, 1, $false, "adequately manage", 2)
$d.Content.Find.Execute("quantative predictions", $true, $false, $false, $false, $false, $true, 1, $false, "quantitative predictions", 2)

# ---------------------------------------------------------------------
# 2. Insert a new sentence (+ its citation) right after "(Hilborn and
#    Walters 1992)" and before the paragraph's closing period.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$end5 = $p5.Range.End
$insertPoint = $d.Range($end5 - 2, $end5 - 2)
$newSentence = ". A wide array data may be collected for an assessment and, formally, an assessment often boils down to algorithms that convert these data to advice for policy decisions. In some cases, particularly for commercially valuable species, this means that hundreds or thousands of historical data points from the monitoring program of a stock gets reduced into a single policy value, such as a recommended catch quota (Maunder, Schnute, and Ianelli 2009)"
$insertPoint.InsertBefore($newSentence)

# ---------------------------------------------------------------------
# 3. Add a new "Body Text" paragraph right after the FirstParagraph.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs.Item(6)
$newParaRange = $d.Range($p6.Range.Start, $p6.Range.End)
$bodyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">Data from monitoring programs&#8230;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newParaRange.InsertXML($bodyXml)

# ---------------------------------------------------------------------
# 4. Add the new Maunder/Schnute/Ianelli reference to the bibliography.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endLast = $lastPara.Range.End
$biblioInsertPoint = $d.Range($endLast, $endLast)
$biblioXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Bibliography"/></w:pPr><w:r><w:t xml:space="preserve">Maunder, Mark N, Jon T Schnute, and James N Ianelli. 2009. &#8220;Computers in Fisheries Population Dynamics.&#8221; In</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Computers in Fisheries Research</w:t></w:r><w:r><w:t xml:space="preserve">, 337&#8211;72. Springer.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$biblioInsertPoint.InsertXML($biblioXml)

# ---------------------------------------------------------------------
# Summary for debugging.
# ---------------------------------------------------------------------
Write-Output "Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "$i : [$($p.Style.NameLocal)] $($p.Range.Text)"
}
